# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" worksheet (with fund-holding data) right after
# "总计" and before "2022-Q3", and adds a corresponding summary row to the
# "总计" (总计) sheet. All other existing quarter sheets keep their own
# data untouched; they simply shift one tab-position to the right, which
# Excel/OOXML does automatically as a side effect of inserting the new
# sheet before them.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# Helper: write a value as genuine text (not a number), without leaving
# any stray cell-style behind (the quote-prefix trick tags the cell with
# a "quotePrefix" style, so we strip formatting right back off again).
function Set-TextCell($sheet, $row, $col, $val) {
    $sheet.Cells.Item($row, $col).Value = "'" + $val
    $sheet.Cells.Item($row, $col).ClearFormats()
}

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet.
#    Use the "2022-Q2" sheet (identical 6-row layout/formatting) as a
#    formatting template, copying it in immediately before "2022-Q3" so
#    tab order becomes: 总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3
# ---------------------------------------------------------------------
$template = $sheets.Item("2022-Q2")
$template.Copy($sheets.Item("2022-Q3"))
$newSheet = $sheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q4"

# Row 2
Set-TextCell $newSheet 2 2 "161611"
Set-TextCell $newSheet 2 3 "融通内需驱动混合A/B"
Set-TextCell $newSheet 2 4 "9.28"
Set-TextCell $newSheet 2 5 "92.61"
Set-TextCell $newSheet 2 6 "3.80"
Set-TextCell $newSheet 2 7 "0.3526"
$newSheet.Cells.Item(2, 8).Value = 9

# Row 3
Set-TextCell $newSheet 3 2 "014109"
Set-TextCell $newSheet 3 3 "融通内需驱动混合C"
Set-TextCell $newSheet 3 4 "3.63"
Set-TextCell $newSheet 3 5 "92.61"
Set-TextCell $newSheet 3 6 "3.80"
Set-TextCell $newSheet 3 7 "0.1379"
$newSheet.Cells.Item(3, 8).Value = 9

# Row 4
Set-TextCell $newSheet 4 2 "014106"
Set-TextCell $newSheet 4 3 "融通成长30灵活配置混合C"
Set-TextCell $newSheet 4 4 "2.19"
Set-TextCell $newSheet 4 5 "94.02"
Set-TextCell $newSheet 4 6 "3.26"
Set-TextCell $newSheet 4 7 "0.0714"
$newSheet.Cells.Item(4, 8).Value = 9

# Row 5
Set-TextCell $newSheet 5 2 "002252"
Set-TextCell $newSheet 5 3 "融通成长30灵活配置混合A/B"
Set-TextCell $newSheet 5 4 "1.65"
Set-TextCell $newSheet 5 5 "94.02"
Set-TextCell $newSheet 5 6 "3.26"
Set-TextCell $newSheet 5 7 "0.0538"
$newSheet.Cells.Item(5, 8).Value = 9

# Row 6
Set-TextCell $newSheet 6 2 "620001"
Set-TextCell $newSheet 6 3 "金元顺安宝石动力混合"
Set-TextCell $newSheet 6 4 "0.46"
Set-TextCell $newSheet 6 5 "56.89"
Set-TextCell $newSheet 6 6 "8.28"
Set-TextCell $newSheet 6 7 "0.0381"
$newSheet.Cells.Item(6, 8).Value = 2

# ---------------------------------------------------------------------
# 2. Add the "2022-Q4" row to the "总计" (summary) sheet, as new row 2,
#    pushing the existing rows down.
# ---------------------------------------------------------------------
$summary = $sheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 5
$summary.Cells.Item(2, 4).Value = 0.65

# The inserted row picked up formatting from the row above (the header);
# reset it to Normal and then copy the correct numeric-index cell style
# from the row below (matches every other data row in the sheet).
$summary.Range("A2:D2").Style = "Normal"
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)

# Renumber the index column (A) for the rows that shifted down, since
# Insert() kept their old values instead of bumping them.
for ($r = 3; $r -le 7; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}
